$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "43.134.89"
Set-TextValue "E2" "  +1.28%  "
Set-TextValue "D3" "2.382.09"
Set-TextValue "E3" "  +6.85%  "
Set-TextValue "E4" "  -0.54%  "
Set-TextValue "D5" "323.75"
Set-TextValue "E5" "  +10.27%  "
Set-TextValue "D6" "105.84"
Set-TextValue "E6" "  -5.59%  "
Set-TextValue "E7" "  +2.98%  "
Set-TextValue "E8" "  +0.01%  "
Set-TextValue "D9" "0.665"
Set-TextValue "E9" "  +11.11%  "
Set-TextValue "D10" "41.57"
Set-TextValue "E10" "  -4.35%  "
Set-TextValue "E11" "  +1.96%  "
Set-TextValue "D12" "8.58"
Set-TextValue "E12" "  -1.45%  "
Set-TextValue "D13" "1.03"
Set-TextValue "E13" "  -2.68%  "
Set-TextValue "D14" "16.92"
Set-TextValue "E14" "  +13.05%  "
Set-TextValue "E15" "  +2.12%  "
Set-TextValue "D16" "2.740.43"
Set-TextValue "E16" "  +6.85%  "
Set-TextValue "D17" "2.372.72"
Set-TextValue "E17" "  +5.75%  "
Set-TextValue "D18" "43.079.86"
Set-TextValue "E18" "  +1.45%  "
Set-TextValue "D19" "7.96"
Set-TextValue "E19" "  +10.17%  "
Set-TextValue "E20" "  +2.69%  "
Set-TextValue "D21" "76.44"
Set-TextValue "E21" "  +3.84%  "
Set-TextValue "D22" "275.48"
Set-TextValue "E22" "  +14.62%  "
Set-TextValue "D23" "3.39"
Set-TextValue "E23" "  +1.33%  "
Set-TextValue "D24" "2.41"
Set-TextValue "E24" "  +2.28%  "
Set-TextValue "D25" "9.49"
Set-TextValue "E25" "  +7.20%  "
Set-TextValue "D26" "11.70"
Set-TextValue "E26" "  +2.31%  "
Set-TextValue "E27" "  +0.03%  "
Set-TextValue "E28" "  +7.12%  "
Set-TextValue "D29" "37.92"
Set-TextValue "E29" "  +2.38%  "
Set-TextValue "D30" "175.13"
Set-TextValue "E30" "  -0.07%  "
Set-TextValue "E31" "  -1.41%  "
Set-TextValue "E32" "  +0.57%  "
Set-TextValue "D33" "0.0921"
Set-TextValue "E33" "  +4.74%  "
Set-TextValue "D34" "5.85"
Set-TextValue "E34" "  +3.15%  "
Set-TextValue "E35" "  +5.58%  "
Set-TextValue "D36" "4.86"
Set-TextValue "E36" "  -1.08%  "
Set-TextValue "D37" "4.20"
Set-TextValue "E37" "  +0.87%  "
Set-TextValue "D38" "0.0366"
Set-TextValue "E38" "  -1.78%  "
Set-TextValue "D39" "0.107"
Set-TextValue "E39" "  +2.76%  "
Set-TextValue "E40" "  +18.68%  "
Set-TextValue "D41" "1.60"
Set-TextValue "E41" "  +22.57%  "
Set-TextValue "D42" "0.231"
Set-TextValue "E42" "  +1.70%  "
Set-TextValue "B43" "MultiversX"
Set-TextValue "C43" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D43" "69.69"
Set-TextValue "E43" "  -2.09%  "
Set-TextValue "B44" "Aave"
Set-TextValue "C44" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "122.97"
Set-TextValue "E44" "  +20.44%  "
Set-TextValue "D45" "94.05"
Set-TextValue "E45" "  +60.72%  "
Set-TextValue "E46" "  +0.07%  "
Set-TextValue "D47" "12.46"
Set-TextValue "E47" "  +1.30%  "
Set-TextValue "E48" "  +11.18%  "
Set-TextValue "E49" "  +2.48%  "
Set-TextValue "D50" "1.31"
Set-TextValue "E50" "  +2.36%  "
Set-TextValue "D51" "1.608.27"
Set-TextValue "E51" "  +12.67%  "
